$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -6
$ws.Range("F9").Value = 1
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = -3
$ws.Range("F20").Value = -1
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = -1
